# Generate Report for Handback
# Updates the "ddae4916-7eb4-4bbf-861d-476a0f304765" row (row 6) on both the
# zh-cn and de-de worksheets: the handback has now completed, but the
# handed-back file's commit is behind the latest source commit, so the
# Latest Target File / Latest Handback File / Latest Handback DateTime /
# Error Detail columns are populated accordingly.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/be8f5ce606e7c3af50338f7fd00c884eebfc6b43/e2e/ddae4916-7eb4-4bbf-861d-476a0f304765.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61c86ec587a96b1107673ac85d8fa8c87c7aa583/e2e/ddae4916-7eb4-4bbf-861d-476a0f304765.md."

function Update-LocSheet {
    param($SheetName, $HandbackFile, $HandbackDateTime, $HyperlinkTarget)

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the Error Detail column (P) to fit the long message.
    $ws.Columns.Item(16).ColumnWidth = 40

    # I6: Latest Target File -> becomes a hyperlink to the handed-back md file.
    $i6 = $ws.Range("I6")
    $i6.Value = "ddae4916-7eb4-4bbf-861d-476a0f304765.md"
    $ws.Hyperlinks.Add($i6, $HyperlinkTarget, "", "", "ddae4916-7eb4-4bbf-861d-476a0f304765.md")

    # J6: Latest Handback File
    $ws.Range("J6").Value = $HandbackFile

    # K6: Latest Handback DateTime
    $ws.Range("K6").Value = $HandbackDateTime

    # P6: Error Detail
    $ws.Range("P6").Value = $errorDetail
}

Update-LocSheet "zh-cn" `
    "ddae4916-7eb4-4bbf-861d-476a0f304765.adc0c177b9a71170813caf24e764586efb9d76a3.zh-cn.xlf" `
    "2016-10-24 09:33:02" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/63eb60573d6fc6d2da429103d0a0c6de0d8f0f1d/e2e/ddae4916-7eb4-4bbf-861d-476a0f304765.md"

Update-LocSheet "de-de" `
    "ddae4916-7eb4-4bbf-861d-476a0f304765.adc0c177b9a71170813caf24e764586efb9d76a3.de-de.xlf" `
    "2016-10-24 09:33:18" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/63eb60573d6fc6d2da429103d0a0c6de0d8f0f1d/e2e/ddae4916-7eb4-4bbf-861d-476a0f304765.md"
